$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1086
$ws.Range("F5").Value = 510
$ws.Range("F6").Value = 217
$ws.Range("F7").Value = 658
$ws.Range("F8").Value = 236
$ws.Range("F10").Value = 77
$ws.Range("F11").Value = 213
$ws.Range("F12").Value = 144
$ws.Range("F13").Value = 1775
$ws.Range("F14").Value = 425
$ws.Range("F15").Value = 39
$ws.Range("F16").Value = 487
$ws.Range("F17").Value = 253
$ws.Range("F18").Value = 408
$ws.Range("F26").Value = 1532
$ws.Range("F27").Value = 273

# Sheet: 演出 (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 206
$ws.Range("F7").Value = 9

# Sheet: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 393

# Sheet: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 393
$ws.Range("F5").Value = 1086
$ws.Range("F8").Value = 510
$ws.Range("F9").Value = 217
$ws.Range("F10").Value = 658
$ws.Range("F12").Value = 236
$ws.Range("F14").Value = 77
$ws.Range("F15").Value = 213
$ws.Range("F16").Value = 144
$ws.Range("F17").Value = 1775
$ws.Range("F18").Value = 206
$ws.Range("F19").Value = 425
$ws.Range("F20").Value = 39
$ws.Range("F21").Value = 487
$ws.Range("F22").Value = 253
$ws.Range("F23").Value = 408
$ws.Range("F27").Value = 9
$ws.Range("F39").Value = 1532
$ws.Range("F40").Value = 273
